$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 3-5: results were re-ordered (same matches, new row order) ---
# Row 3
$ws.Range("F3").Value = "Persepolis"
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = "Aluminium Arak"
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 1.37
$ws.Range("K3").Value = "02/08/2023 06:12"
$ws.Range("L3").Value = 1.46
$ws.Range("M3").Value = "09/08/2023 16:33"
$ws.Range("N3").Value = 3.83
$ws.Range("O3").Value = "02/08/2023 06:12"
$ws.Range("P3").Value = 3.63
$ws.Range("Q3").Value = "09/08/2023 16:33"
$ws.Range("R3").Value = 8.02
$ws.Range("S3").Value = "02/08/2023 06:12"
$ws.Range("T3").Value = 9.26
$ws.Range("U3").Value = "09/08/2023 16:33"
$ws.Range("V3").Value = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/persepolis-aluminium-arak/xpKbcLUC/"

# Row 4
$ws.Range("F4").Value = "Malavan"
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = "Zob Ahan"
$ws.Range("I4").Value = 2
$ws.Range("J4").Value = 2.62
$ws.Range("K4").Value = "08/08/2023 06:12"
$ws.Range("L4").Value = 3.15
$ws.Range("M4").Value = "09/08/2023 16:52"
$ws.Range("N4").Value = 2.62
$ws.Range("O4").Value = "08/08/2023 06:12"
$ws.Range("P4").Value = 2.43
$ws.Range("Q4").Value = "09/08/2023 16:52"
$ws.Range("R4").Value = 2.86
$ws.Range("S4").Value = "08/08/2023 06:12"
$ws.Range("T4").Value = 2.92
$ws.Range("U4").Value = "09/08/2023 16:52"
$ws.Range("V4").Value = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/malavan-zob-ahan/4din2qGs/"

# Row 5
$ws.Range("F5").Value = "Gol Gohar"
$ws.Range("G5").Value = 3
$ws.Range("H5").Value = "Foolad"
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 2.62
$ws.Range("K5").Value = "08/08/2023 06:12"
$ws.Range("L5").Value = 2.77
$ws.Range("M5").Value = "09/08/2023 17:47"
$ws.Range("N5").Value = 2.63
$ws.Range("O5").Value = "08/08/2023 06:12"
$ws.Range("P5").Value = 2.54
$ws.Range("Q5").Value = "09/08/2023 17:47"
$ws.Range("R5").Value = 2.85
$ws.Range("S5").Value = "08/08/2023 06:12"
$ws.Range("T5").Value = 3.17
$ws.Range("U5").Value = "09/08/2023 17:47"
$ws.Range("V5").Value = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/gol-gohar-foolad/xAkf0Npf/"

# Row 35
$ws.Range("F35").Value = "Shams Azar Qazvin"
$ws.Range("G35").Value = 1
$ws.Range("H35").Value = "Zob Ahan"
$ws.Range("I35").Value = 2
$ws.Range("J35").Value = 2.98
$ws.Range("K35").Value = "05/10/2023 09:42"
$ws.Range("L35").Value = 2.83
$ws.Range("M35").Value = "05/10/2023 15:32"
$ws.Range("N35").Value = 2.78
$ws.Range("O35").Value = "05/10/2023 09:42"
$ws.Range("P35").Value = 2.65
$ws.Range("Q35").Value = "05/10/2023 15:16"
$ws.Range("R35").Value = 2.61
$ws.Range("S35").Value = "05/10/2023 09:42"
$ws.Range("T35").Value = 2.93
$ws.Range("U35").Value = "05/10/2023 15:32"
$ws.Range("V35").Value = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/shams-azar-qazvin-zob-ahan/Aq7VT9AN/"

# Row 36
$ws.Range("F36").Value = "Aluminium Arak"
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = "Tractor"
$ws.Range("I36").Value = 4
$ws.Range("J36").Value = 2.83
$ws.Range("K36").Value = "05/10/2023 09:42"
$ws.Range("L36").Value = 3.18
$ws.Range("M36").Value = "05/10/2023 15:43"
$ws.Range("N36").Value = 2.74
$ws.Range("O36").Value = "05/10/2023 09:42"
$ws.Range("P36").Value = 2.45
$ws.Range("Q36").Value = "05/10/2023 15:43"
$ws.Range("R36").Value = 2.71
$ws.Range("S36").Value = "05/10/2023 09:42"
$ws.Range("T36").Value = 2.87
$ws.Range("U36").Value = "05/10/2023 15:43"
$ws.Range("V36").Value = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/aluminium-arak-tractor/466ZSTPT/"

# Row 62
$ws.Range("F62").Value = "Mes Rafsanjan"
$ws.Range("G62").Value = 3
$ws.Range("H62").Value = "Foolad"
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 2.2
$ws.Range("K62").Value = "02/11/2023 00:42"
$ws.Range("L62").Value = 2.3
$ws.Range("M62").Value = "03/11/2023 12:29"
$ws.Range("N62").Value = 2.65
$ws.Range("O62").Value = "02/11/2023 00:42"
$ws.Range("P62").Value = 2.48
$ws.Range("Q62").Value = "03/11/2023 12:29"
$ws.Range("R62").Value = 3.55
$ws.Range("S62").Value = "02/11/2023 00:42"
$ws.Range("T62").Value = 4.29
$ws.Range("U62").Value = "03/11/2023 12:29"
$ws.Range("V62").Value = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/mes-rafsanjan-foolad/z7Iven9M/"

# Row 63
$ws.Range("F63").Value = "Havadar SC"
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = "Paykan"
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 2.17
$ws.Range("K63").Value = "02/11/2023 00:42"
$ws.Range("L63").Value = 2.15
$ws.Range("M63").Value = "03/11/2023 12:28"
$ws.Range("N63").Value = 2.65
$ws.Range("O63").Value = "02/11/2023 00:42"
$ws.Range("P63").Value = 2.41
$ws.Range("Q63").Value = "03/11/2023 12:28"
$ws.Range("R63").Value = 3.62
$ws.Range("S63").Value = "02/11/2023 00:42"
$ws.Range("T63").Value = 4.31
$ws.Range("U63").Value = "03/11/2023 12:28"
$ws.Range("V63").Value = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/havadar-sc-paykan/KvfBkQWq/"

# --- Append two new match rows (71, 72); copy formatting from the last existing row first ---
$ws.Range("A70:V70").Copy()
$ws.Range("A71:V72").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 71
$ws.Range("A71").Value = 70
$ws.Range("B71").Value = "iran"
$ws.Range("C71").Value = "persian-gulf-pro-league"
$ws.Range("D71").Value = "2023-2024"
$ws.Range("E71").Value = 45242.52083333334
$ws.Range("F71").Value = "Nassaji Mazandaran"
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = "Zob Ahan"
$ws.Range("I71").Value = 2
$ws.Range("J71").Value = 2.52
$ws.Range("K71").Value = "11/11/2023 00:43"
$ws.Range("L71").Value = 3.23
$ws.Range("M71").Value = "12/11/2023 12:26"
$ws.Range("N71").Value = 2.63
$ws.Range("O71").Value = "11/11/2023 00:43"
$ws.Range("P71").Value = 2.38
$ws.Range("Q71").Value = "12/11/2023 12:26"
$ws.Range("R71").Value = 3.06
$ws.Range("S71").Value = "11/11/2023 00:43"
$ws.Range("T71").Value = 2.92
$ws.Range("U71").Value = "12/11/2023 12:26"
$ws.Range("V71").Value = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/mazandaran-zob-ahan/Q7ibU3xj/"

# Row 72
$ws.Range("A72").Value = 71
$ws.Range("B72").Value = "iran"
$ws.Range("C72").Value = "persian-gulf-pro-league"
$ws.Range("D72").Value = "2023-2024"
$ws.Range("E72").Value = 45242.58333333334
$ws.Range("F72").Value = "Sepahan"
$ws.Range("G72").Value = 1
$ws.Range("H72").Value = "Persepolis"
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 2.31
$ws.Range("K72").Value = "11/11/2023 02:12"
$ws.Range("L72").Value = 2.34
$ws.Range("M72").Value = "12/11/2023 13:57"
$ws.Range("N72").Value = 2.8
$ws.Range("O72").Value = "11/11/2023 02:12"
$ws.Range("P72").Value = 2.8
$ws.Range("Q72").Value = "12/11/2023 13:57"
$ws.Range("R72").Value = 3.1
$ws.Range("S72").Value = "11/11/2023 02:12"
$ws.Range("T72").Value = 3.48
$ws.Range("U72").Value = "12/11/2023 13:57"
$ws.Range("V72").Value = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/sepahan-persepolis/vBe2TNid/"

